$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Formatting cleanup -------------------------------------------------
# A handful of cells carried a stray "Serif" / blank font instead of the
# sheet's normal Arial formatting. Bring them in line with the rest of the
# table before adding the new data row.

# B6, B8, D8 were using a blank/default font -> match the normal cell font
# (e.g. A6, which already carries the standard formatting).
$ws.Range("A6").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("A6").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("A6").Copy()
$ws.Range("D8").PasteSpecial(-4122)

# F6, G6, F7, G7, F36, G36, F37, G37 were Serif without a color -> align
# them with H28's Serif + automatic color formatting.
$ws.Range("H28").Copy()
$ws.Range("F6").PasteSpecial(-4122)
$ws.Range("H28").Copy()
$ws.Range("G6").PasteSpecial(-4122)
$ws.Range("H28").Copy()
$ws.Range("F7").PasteSpecial(-4122)
$ws.Range("H28").Copy()
$ws.Range("G7").PasteSpecial(-4122)
$ws.Range("H28").Copy()
$ws.Range("F36").PasteSpecial(-4122)
$ws.Range("H28").Copy()
$ws.Range("G36").PasteSpecial(-4122)
$ws.Range("H28").Copy()
$ws.Range("F37").PasteSpecial(-4122)
$ws.Range("H28").Copy()
$ws.Range("G37").PasteSpecial(-4122)

# --- New data --------------------------------------------------------
# Added new Earth data! (Day 3, 130mm/s)
$ws.Range("A9:G9").Copy()
$ws.Range("A38:G38").PasteSpecial(-4122)

$ws.Range("A38").Value = "N/A"
$ws.Range("B38").Value = "Day3-Earth-130mms.mov"
$ws.Range("C38").Value = 3
$ws.Range("D38").Value = "Earth"
$ws.Range("E38").Value = 130
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 7
